# Reposition/resize the pictures on slides 1 and 2.
# Left/Top/Width/Height are expressed in points in the PowerPoint object
# model (1 pt = 12700 EMU). The literal values below are chosen so that,
# after the host's internal float32 storage, they convert back to the
# exact target EMU values used by the canonical OOXML.

$p = $ppt.ActivePresentation

# Slide 1
$s1 = $p.Slides.Item(1)

$pic5 = $s1.Shapes.Item("Picture 5")
$pic5.Left   = 31.970866241732285
$pic5.Top    = 73.12945181889765
$pic5.Width  = 434.2367716535433
$pic5.Height = 235.77937007874016

$pic19 = $s1.Shapes.Item("Picture 19")
$pic19.Left   = 31.970944881889764
$pic19.Top    = 310.2585039370079
$pic19.Width  = 434.2367716535433
$pic19.Height = 229.74149606299213

# Slide 2
$s2 = $p.Slides.Item(2)

$pic7 = $s2.Shapes.Item("Picture 7")
$pic7.Left   = 10.832283564566929
$pic7.Top    = 58.38708691417323
$pic7.Width  = 444.72622047244096
$pic7.Height = 248.175911011811

$pic14 = $s2.Shapes.Item("Picture 14")
$pic14.Left   = 10.832283564566929
$pic14.Top    = 312.9788976377953
$pic14.Width  = 329.0463105125984
$pic14.Height = 232.36771653543306
